$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games), rows 2-8, columns B-G
$data = @{
    2 = @(3.230985683306322, 10.29869402782916, 3.900430680208489, 8.660232485948974, 1, 26.09034287729295)
    3 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
    4 = @(0.3048080303191223, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 2.626907116734944)
    5 = @(0.04763786555579896, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0, 1.014732764554632)
    6 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 3.781711156805759)
    7 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 3.781711156805759)
    8 = @(0.127881588408715, 0.3127903958511391, 26.21740644021617, 8.660232485948974, 0, 35.318310910425)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
